# Applies the "automatic update" diff to the Avverkningsanmälningar sheet:
#  - column C (Förändrad) bumps from 46062 to 46063 for every data row (2-14)
#  - rows 6-14 (except row 7) get reshuffled to new A/B/G values
#  - two new rows (15, 16) are appended with fresh case data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style/number-format reference cells (existing date + wrap-text columns)
$dateFormat = $ws.Range("B2").NumberFormat

function Set-Row($r, $a, $b, $g) {
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = 46063
    $ws.Range("G$r").Value = $g
}

# Rows 2-5 and 7 keep their own Beteckning/Datum/Area; only column C (Förändrad) changes
foreach ($r in @(2,3,4,5,7)) {
    $ws.Range("C$r").Value = 46063
}

# Rows 6, 8-14 are reshuffled with new Beteckning/Datum/Area values
Set-Row 6 "A 5792-2024" 45335 5.6
Set-Row 8 "A 2593-2024" 45313.69204861111 2.3
Set-Row 9 "A 7333-2025" 45703.35899305555 0.9
Set-Row 10 "A 35642-2023" 45147 1.2
Set-Row 11 "A 28288-2023" 45099.6349537037 0.5
Set-Row 12 "A 8194-2025" 45708 1.9
Set-Row 13 "A 50997-2025" 45946 1.5
Set-Row 14 "A 13651-2023" 45006 2.2

# Row 14 gains an explicit row height (matches the other data rows now)
$ws.Rows.Item(14).RowHeight = 15

# New rows 15 and 16
Set-Row 15 "A 7814-2026" 46062.61388888889 1.1
Set-Row 16 "A 7827-2026" 46062.63958333333 2.1

foreach ($r in 15..16) {
    $ws.Range("D$r").Value = "SKÅNE LÄN"
    $ws.Range("E$r").Value = "HELSINGBORG"
    foreach ($col in @("H","I","J","K","L","M","N","O","P","Q")) {
        $ws.Range("$col$r").Value = 0
    }
    $ws.Range("B$r").NumberFormat = $dateFormat
    $ws.Range("C$r").NumberFormat = $dateFormat
    $ws.Range("R$r").Value = ""
    $ws.Range("R$r").WrapText = $true
}

# Only row 15 gets an explicit custom row height (row 16 stays on the sheet default)
$ws.Rows.Item(15).RowHeight = 15
